$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark from its current location
#    (between "Back to" and " sentiment analysis").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Append two new paragraphs at the very end of the document body:
#      - an empty paragraph
#      - a paragraph containing "chANGEEEESSSS", with the _GoBack bookmark
#        placed right after the text (a zero-length bookmark at the end
#        of the paragraph).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$emptyPara = $d.Paragraphs.Last
$emptyPara.Range.InsertParagraphAfter()

$finalPara = $d.Paragraphs.Last
# Insert the visible text plus a throwaway sentinel character so that the
# bookmark can be anchored strictly *inside* a run (anchoring a zero-length
# bookmark exactly at a run/paragraph boundary is unreliable), then remove
# the sentinel afterwards.
$finalPara.Range.InsertAfter("chANGEEEESSSS" + "@@SENTINEL@@")

$p = $d.Paragraphs.Last
$pos = $p.Range.End - (("@@SENTINEL@@").Length + 1)
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$findRange = $d.Content
[void]$findRange.Find.Execute("@@SENTINEL@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Delete()
